# Add a new daily-report sheet "2019-7-21" by duplicating "2019-7-20"
# and filling in the new day's content.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("2019-7-20")
$src.Copy($null, $src)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2019-7-21"

# Header date fields (plan date / actual date / the "date:" row).
$newSheet.Range("F3").Value = "2019.07.21"
$newSheet.Range("D6").Value = "2019.07.21"
$newSheet.Range("D7").Value = "2019.07.21"

# Row 10: 学习有关软件架构的相关知识
$newSheet.Range("D10").Value = "学习有关软件架构的相关知识"
$newSheet.Range("G10").Value = 3
$newSheet.Range("H10").Value = "学习和习题联系"
$newSheet.Range("K10").Value = "是"

# Row 11: 编写代码，英文前缀树存储
$newSheet.Range("D11").Value = "编写代码，英文前缀树存储"
$newSheet.Range("G11").Value = 3
$newSheet.Range("H11").Value = "编写代码"
$newSheet.Range("K11").Value = "是"

# Row 12: 编写中文对照表
$newSheet.Range("D12").Value = "编写中文对照表"
$newSheet.Range("G12").Value = 0.5
$newSheet.Range("H12").Value = "编写代码"
$newSheet.Range("K12").Value = "否"

# Selection moves to N14 on the new sheet.
$newSheet.Range("N14").Select()
